# Add RecipientSignature and RecipientPhoto columns on the "mark as
# delivered" import template, plus a wrap-text styled placeholder row
# underneath the header so users have somewhere to paste a signature /
# photo reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (adds two shared strings: RecipientSignature, RecipientPhoto)
$ws.Range("F1").Value = "RecipientSignature"
$ws.Range("G1").Value = "RecipientPhoto"

# Widen the columns to fit the new, longer headers
$ws.Columns("A:E").ColumnWidth = 20.77734375
$ws.Columns("F:G").ColumnWidth = 25.77734375

# Second row, first cell gets a dedicated wrap-text style (9pt grey Calibri)
$cell = $ws.Range("A2")
$cell.WrapText = $true
$cell.Font.Size = 9
$cell.Font.Color = 3355443

# Leave the selection where the author left it
$ws.Range("F8").Select() | Out-Null
